$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.086084842681885
$ws.Range("B1").Value = 6.08387565612793
$ws.Range("C1").Value = 3.139637231826782
$ws.Range("D1").Value = 1.383564710617065
$ws.Range("E1").Value = 0.9713207483291626
